$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Target values reconstructed from the diff for changed rows (row 34, and rows 175-253)
# Columns: A = Codigo IE (number), B = Nivel Plan, C = Nivel Lider, D = Nivel Directivos, E = Dimension 5
$rows = @(
    @(34, 33, "1A", "1A", "1A", "1A"),
    @(175, 174, "1A", "1A", "5", "1A"),
    @(176, 175, "1B", "1B", "5", "1B"),
    @(177, 176, "1A", "1A", "1A", "1A"),
    @(178, 177, "1A", "1A", "5", "1A"),
    @(179, 178, "1A", "1A", "5", "1A"),
    @(180, 179, "1A", "1A", "5", "1A"),
    @(181, 180, "1A", "1A", "5", "1A"),
    @(182, 181, "1A", "1A", "1A", "1A"),
    @(183, 182, "1A", "1A", "5", "1A"),
    @(184, 183, "1A", "1A", "5", "1A"),
    @(185, 184, "1A", "1A", "5", "1A"),
    @(186, 185, "1A", "1A", "5", "1A"),
    @(187, 186, "1A", "5", "5", "1A"),
    @(188, 187, "1A", "1A", "5", "1A"),
    @(189, 188, "1A", "1A", "1A", "1A"),
    @(190, 189, "1A", "1A", "5", "1A"),
    @(191, 190, "1A", "1A", "1A", "1A"),
    @(192, 191, "1B", "1A", "1A", "1A"),
    @(193, 192, "1A", "1A", "5", "1A"),
    @(194, 193, "1A", "1A", "5", "1A"),
    @(195, 194, "1A", "1A", "1A", "1A"),
    @(196, 195, "1A", "1A", "5", "1A"),
    @(197, 196, "1A", "1A", "1A", "1A"),
    @(198, 197, "1A", "1A", "1A", "1A"),
    @(199, 198, "1A", "1A", "5", "1A"),
    @(200, 199, "1A", "1A", "5", "1A"),
    @(201, 200, "1A", "1A", "1A", "1A"),
    @(202, 201, "1A", "1A", "1A", "1A"),
    @(203, 202, "1A", "1A", "5", "1A"),
    @(204, 203, "1A", "1A", "5", "1A"),
    @(205, 204, "1A", "1A", "5", "1A"),
    @(206, 205, "1A", "1A", "5", "1A"),
    @(207, 206, "1A", "1A", "1A", "1A"),
    @(208, 207, "1A", "1A", "1A", "1A"),
    @(209, 208, "1A", "1A", "5", "1A"),
    @(210, 209, "1A", "1A", "5", "1A"),
    @(211, 210, "1A", "1A", "1A", "1A"),
    @(212, 211, "1A", "1A", "5", "1A"),
    @(213, 212, "1A", "1A", "5", "1A"),
    @(214, 213, "1B", "1A", "5", "1A"),
    @(215, 214, "1A", "1A", "1A", "1A"),
    @(216, 215, "1A", "1A", "5", "1A"),
    @(217, 216, "1A", "1A", "5", "1A"),
    @(218, 217, "1A", "1A", "5", "1A"),
    @(219, 218, "1A", "1A", "1A", "1A"),
    @(220, 219, "1A", "1A", "1A", "1A"),
    @(221, 220, "1A", "1A", "5", "1A"),
    @(222, 221, "1A", "1A", "5", "1A"),
    @(223, 222, "1B", "1B", "5", "1B"),
    @(224, 223, "1A", "1A", "1A", "1A"),
    @(225, 224, "1B", "1A", "5", "1A"),
    @(226, 225, "1A", "1A", "5", "1A"),
    @(227, 226, "1A", "5", "5", "1A"),
    @(228, 227, "1A", "1A", "5", "1A"),
    @(229, 228, "1A", "1A", "5", "1A"),
    @(230, 229, "1A", "1B", "5", "1A"),
    @(231, 230, "1A", "1A", "1A", "1A"),
    @(232, 231, "1A", "1A", "5", "1A"),
    @(233, 232, "1A", "1A", "5", "1A"),
    @(234, 233, "1A", "1A", "5", "1A"),
    @(235, 234, "1A", "1A", "1A", "1A"),
    @(236, 235, "1A", "1A", "5", "1A"),
    @(237, 236, "1A", "1A", "5", "1A"),
    @(238, 237, "1A", "1A", "5", "1A"),
    @(239, 238, "1A", "1A", "5", "1A"),
    @(240, 239, "1A", "1A", "1A", "1A"),
    @(241, 240, "5", "5", "1A", "1A"),
    @(242, 241, "1A", "1A", "1A", "1A"),
    @(243, 242, "4", "1A", "5", "1A"),
    @(244, 243, "5", "1B", "5", "1B"),
    @(245, 244, "5", "1B", "5", "1B"),
    @(246, 245, "1A", "1A", "5", "1A"),
    @(247, 246, "1A", "1A", "1A", "1A"),
    @(248, 247, "1A", "1A", "1A", "1A"),
    @(249, 248, "1A", "1A", "5", "1A"),
    @(250, 249, "1A", "1A", "5", "1A"),
    @(251, 250, "1A", "1A", "5", "1A"),
    @(252, 251, "1A", "1A", "5", "1A"),
    @(253, 252, "1A", "1A", "5", "1A")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
}
